$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Update the text of RQ8, RQ9 and RQ10 (rows 10, 11, 12 / B10:B12) ---
$ws.Range("B10").Value = "RQ8: Si el operario desea editar su información personal, el sistema debe permitirlo, a excepción de su C.I/C.C"
$ws.Range("B11").Value = "RQ9: Si el cliente desea cambiar algun dato personal registrado en el sistema, el sistema debe comunicar al operario el cual realizará el cambio de la información en el sistema."
$ws.Range("B12").Value = "RQ10: El sistema debe realizar el cálculo de arqueo de caja en base a la información ingresada por el usuario de la cantidad de billetes y monedas recibidas. "

# --- Row height adjustments caused by the new (longer/shorter) wrapped text ---
$ws.Rows.Item(6).RowHeight = 75
$ws.Rows.Item(10).RowHeight = 45
$ws.Rows.Item(11).RowHeight = 94.5
$ws.Rows.Item(12).RowHeight = 69

# --- Keep the B12 cell selected/active, matching the saved view state ---
$ws.Range("B12").Select()
